$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps (rows 36/37 and 42/43) ---
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# --- Price (D) and Volume(1h) (E) updates ---
# Force Text number format first so Excel does not reinterpret values such as
# "294.27" or "0.0779" as real numbers (the source data stores them as text).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.623.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.417.15"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "294.27"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -5.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.51%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.76"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0779"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.05"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.791.58"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.414.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.839"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.547.42"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.28"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.08"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.50"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.81"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -6.89%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +14.27%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.63"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.99"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +17.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.114"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.116"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.57"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.76"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.987.39"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "88.83"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.09"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +26.43%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "100.56"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.660.52"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.98%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.39%  "
